# Updates cryptocurrency price/volume figures on the "cryptos" sheet
# (refreshed symbol-list snapshot). Values are written as text, matching
# the original inline-string cell storage (e.g. "-0.68%" stays literal text,
# not a numeric percentage), and the default "Normal" style is restored so
# no stray number formatting/quote-prefix style is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    # Leading apostrophe forces Excel to store the value as literal text
    # instead of re-interpreting "44.39" or "-0.68%" as a number/percentage.
    $cell.Formula = "'" + $value
    # Reset to the default style so the quote-prefix text entry does not
    # leave behind a different cell style than the other untouched cells.
    $cell.Style = "Normal"
}

Set-TextValue "D2" "328.66"
Set-TextValue "E2" "-0.68%"
Set-TextValue "D3" "44.39"
Set-TextValue "E3" "6.98%"
Set-TextValue "D4" "5.480"
Set-TextValue "E4" "-3.17%"
Set-TextValue "D5" "0.08167"
Set-TextValue "E5" "-2.06%"
Set-TextValue "D7" "4.324"
Set-TextValue "E7" "-3.96%"
Set-TextValue "D8" "1.907"
Set-TextValue "E8" "-4.78%"
Set-TextValue "E9" "-2.92%"
Set-TextValue "D10" "0.9414"
Set-TextValue "E10" "1.64%"
Set-TextValue "D11" "0.1195"
Set-TextValue "E11" "-7.64%"
Set-TextValue "E12" "-3.60%"
Set-TextValue "D13" "0.09874"
Set-TextValue "E13" "4.94%"
Set-TextValue "D14" "0.04181"
Set-TextValue "E14" "7.15%"
Set-TextValue "D15" "0.1068"
Set-TextValue "E15" "1.04%"
Set-TextValue "D16" "0.001290"
Set-TextValue "E16" "-1.02%"
Set-TextValue "D17" "0.006031"
Set-TextValue "E17" "-0.43%"
Set-TextValue "D18" "3.536"
Set-TextValue "E18" "2.80%"
Set-TextValue "D19" "0.3502"
Set-TextValue "E19" "-1.03%"
Set-TextValue "D20" "8.790"
Set-TextValue "E20" "2.86%"
Set-TextValue "D21" "0.1351"
Set-TextValue "E21" "-0.21%"
Set-TextValue "D22" "0.2496"
Set-TextValue "E22" "2.17%"
Set-TextValue "D23" "0.04379"
Set-TextValue "E23" "-0.83%"
Set-TextValue "D24" "0.001240"
Set-TextValue "E24" "-2.77%"
Set-TextValue "D25" "0.004312"
Set-TextValue "D26" "0.0001236"
Set-TextValue "D27" "0.0004005"
Set-TextValue "E27" "31.51%"
Set-TextValue "D39" "0.02713"
Set-TextValue "E39" "-3.02%"
Set-TextValue "D40" "0.05702"
Set-TextValue "E40" "3.03%"
Set-TextValue "D41" "0.007883"
Set-TextValue "E41" "0.96%"
Set-TextValue "D42" "0.009739"
Set-TextValue "E42" "4.61%"
Set-TextValue "D43" "0.1413"
Set-TextValue "E43" "-1.67%"
Set-TextValue "D44" "0.002105"
Set-TextValue "E44" "1.64%"
Set-TextValue "D45" "0.009676"
Set-TextValue "E45" "-12.67%"
Set-TextValue "D46" "0.00007122"
Set-TextValue "E46" "1.75%"
Set-TextValue "E47" "0.38%"
Set-TextValue "D48" "0.003444"
Set-TextValue "E48" "5.28%"
Set-TextValue "D49" "0.002279"
Set-TextValue "E49" "-0.05%"
Set-TextValue "D50" "0.00002110"
Set-TextValue "E50" "0.38%"
Set-TextValue "D51" "0.0002010"
Set-TextValue "E51" "0.38%"
